$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, border, centered) from H1 to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Data values for the new columns
$valuesI = @(6, 4, 9, 3, 5, 1, 1, 4)
$valuesJ = @(7, 5, 9, 7, 8, 4, 3, 5)

for ($i = 0; $i -lt $valuesI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $valuesI[$i]
    $ws.Cells.Item($row, 10).Value = $valuesJ[$i]
}
